$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "56.636.83"
$ws.Range("E2").Value = "  +4.28%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.011.10"
$ws.Range("E3").Value = "  +4.66%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"

# Row 5 - BNB
$ws.Range("D5").Value = "508.01"
$ws.Range("E5").Value = "  +8.36%  "

# Row 6 - Solana
$ws.Range("D6").Value = "137.52"
$ws.Range("E6").Value = "  +9.70%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +7.74%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +15.55%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +13.57%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +8.17%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.129"
$ws.Range("E12").Value = "  +5.73%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.521.75"
$ws.Range("E13").Value = "  +4.76%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "25.52"
$ws.Range("E14").Value = "  +10.51%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.0000155"
$ws.Range("E15").Value = "  +16.81%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "56.638.02"
$ws.Range("E16").Value = "  +4.40%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.003.60"
$ws.Range("E17").Value = "  +4.53%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "5.83"
$ws.Range("E18").Value = "  +9.75%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  +10.63%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.86"
$ws.Range("E20").Value = "  +11.77%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "328.40"
$ws.Range("E21").Value = "  +10.96%  "

# Row 22 - Dai
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.06%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.480"
$ws.Range("E23").Value = "  +9.73%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "62.53"
$ws.Range("E24").Value = "  +7.29%  "

# Row 25 - Kaspa
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +12.54%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  -0.06%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0₃0919"
$ws.Range("E27").Value = "  +14.77%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "6.59"
$ws.Range("E28").Value = "  +8.76%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.00"
$ws.Range("E29").Value = "  +14.53%  "

# Row 30 - Fetch.AI
$ws.Range("D30").Value = "1.26"
$ws.Range("E30").Value = "  +13.58%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +10.37%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "20.76"
$ws.Range("E32").Value = "  +11.91%  "

# Row 33 - Monero
$ws.Range("D33").Value = "156.06"
$ws.Range("E33").Value = "  +12.83%  "

# Row 34 - NEARProtocol
$ws.Range("D34").Value = "4.51"
$ws.Range("E34").Value = "  +9.28%  "

# Row 35 - Aptos
$ws.Range("E35").Value = "  +5.11%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +5.30%  "

# Row 37 and 38 swap: Hedera <-> EnergySwap
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "24.29"
$ws.Range("E37").Value = "  +6.29%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.0674"
$ws.Range("E38").Value = "  +10.11%  "

# Row 39 - RenzoRestakedETH
$ws.Range("D39").Value = "3.045.56"
$ws.Range("E39").Value = "  +5.12%  "

# Row 40 - OKB
$ws.Range("D40").Value = "36.67"
$ws.Range("E40").Value = "  +4.86%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.651"
$ws.Range("E42").Value = "  +8.30%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.271.47"
$ws.Range("E43").Value = "  +11.97%  "

# Row 44 - ONDO
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +7.25%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +8.35%  "

# Row 46 - Filecoin
$ws.Range("D46").Value = "3.62"
$ws.Range("E46").Value = "  +7.70%  "

# Row 47 - dogwifhat
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  +24.47%  "

# Row 48 - VeChain
$ws.Range("D48").Value = "0.0237"
$ws.Range("E48").Value = "  +11.25%  "

# Row 49 - Cosmos
$ws.Range("E49").Value = "  +8.67%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "19.17"
$ws.Range("E50").Value = "  +8.64%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +12.02%  "
